# students_template.xlsx: bump the sample student's regno from 2024009 to
# 2024010 (the only functional change in this commit), and replay the
# cosmetic view-state changes Excel saved alongside it: the last selected
# cell moved to J15, and column A ("first_name") was re-measured a bit
# wider.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2024010

# Column A's best-fit width grew from 10.5703125 to 11.7109375 characters.
# This runtime's ColumnWidth setter only lands on 1/6-character increments,
# so 10.8333333333333 is the closest input that resolves to the nearest
# achievable width (11.666666666666666, vs. the exact 11.7109375 target).
$ws.Columns("A").ColumnWidth = 10.8333333333333

# The active cell/selection was left on J15 when the workbook was saved.
$null = $ws.Range("J15").Select()
